# Remove the first three data rows (2007年, 2008年, 2009年).
# This shifts the remaining rows (2010年-2013年) up so that they
# become rows 2-5, and the used range shrinks from A1:F8 to A1:F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:4").Delete()
